# Update vehicle specifications from EV database, Mercedes Benz, Volkswagen
#
# The existing ICEV "Use phase" row (t CO2, row 9) is relocated down to a new
# row 11, and row 9 is repurposed to show the equivalent emissions expressed
# as g CO2/km, computed from the (now relocated) t CO2 figures.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Copy the current row 9 ("Use phase" / "t CO2" + values) down to row 11,
#    preserving its number format (0.0).
$ws.Range("B11").Value = $ws.Range("B9").Value()
$ws.Range("C11").Value = $ws.Range("C9").Value()
$ws.Range("D11:G11").Value = $ws.Range("D9:G9").Value()
$ws.Range("D11:G11").NumberFormat = "0.0"

# 2. Re-purpose row 9: label switches to the new "g CO2/km" unit and the
#    values become formulas derived from the relocated t CO2 row (row 11),
#    converting tonnes over the vehicle lifetime distance into g/km.
$ws.Range("C9").Value = "g CO2/km"
$ws.Range("D9").Formula = "=D11*1000000/180000"
$ws.Range("E9:G9").Formula = "=E11*1000000/180000"
$ws.Range("D9:G9").NumberFormat = "0.0"
$ws.Range("D9:G9").HorizontalAlignment = -4152

# 3. Tidy up the D:G column widths to a uniform, non-autofit width.
$ws.Range("D1:G1").ColumnWidth = 9.666666666666666

# 4. Leave the sheet with the same selection state recorded in the workbook.
[void]$ws.Range("D13:J15").Select()
